$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8ac8d5eb50eda5b218044b7fbfaca67952ba8b3a"

function Update-LangSheet($ws, $xlfFileName, $handoffDateTime) {
    # Status column (B2): "Handoff transform failed" -> "Ready for handoff"
    $ws.Range("B2").Value = "Ready for handoff"

    # New "Latest Handoff File" hyperlink (C2)
    $ws.Hyperlinks.Add($ws.Range("C2"), "$baseUrl/$xlfFileName", "", "", $xlfFileName)

    # Latest Handoff Datetime (D2)
    $ws.Range("D2").Value = $handoffDateTime

    # Handoff Reason (H2): "Ignored" -> "Include"
    $ws.Range("H2").Value = "Include"
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

Update-LangSheet $wsZh "5f07ccd3-7c76-4477-ae83-43e3c1a6f184.47e7d9f23df2305d285788acb1a44c1df18f566d.zh-cn.xlf" "2016-01-15 14:43:59"
Update-LangSheet $wsDe "5f07ccd3-7c76-4477-ae83-43e3c1a6f184.47e7d9f23df2305d285788acb1a44c1df18f566d.de-de.xlf" "2016-01-15 14:44:09"
